$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "Complete" column header, matching the bold style used by the other headers (A3:C3)
$ws.Range("D3").Value = "Complete"
$ws.Range("D3").Font.Bold = $true

# Fix task description typo: "zip code -county" -> "zip code-county"
$ws.Range("A4").Value = "Get zip code-county mapping and write logic to lookup based on report from inpatient data"

# Fill in the paper draft task: owner and completion date
$ws.Range("B12").Value = "Chris"
$ws.Range("C12").Value = Get-Date -Year 2014 -Month 4 -Day 7 -Hour 0 -Minute 0 -Second 0

# Scroll/selection state matching the saved view
$ws.Range("A5").Select()
